# "Cleaning code for v.0.1 part I." -- append two more observer records
# (rows 25 & 26) to Sheet1, just below the existing 23-row table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shared strings get appended to xl/sharedStrings.xml in first-use order,
# and the target workbook expects "RP" before "BO" (indices 26 and 27
# respectively), so write the B column in that order.
$ws.Cells.Item(26, 2).Value = "RP"
$ws.Cells.Item(25, 2).Value = "BO"

$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 3).Value = 2

$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 3).Value = 0

# Match the existing table's row height so the new rows carry the same
# explicit ht="23.25" customHeight="1" as every other data row.
$ws.Rows.Item(25).RowHeight = 23.25
$ws.Rows.Item(26).RowHeight = 23.25

# Scroll the view down and leave the selection on C27 (just past the new
# last row), matching the author's final cursor position.
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("C27").Select()
